$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01177902937386835
$ws.Range("C2").Value = 0.3233677337305382
$ws.Range("D2").Value = 0.1387525778210147
$ws.Range("E2").Value = 0.3724950708680783
$ws.Range("F2").Value = 0.3863630867095333
$ws.Range("G2").Value = 14
$ws.Range("B3").Value = -0.008234150994489252
$ws.Range("C3").Value = 0.3656812685479343
$ws.Range("D3").Value = 0.1736781030421116
$ws.Range("E3").Value = 0.4167470492302394
$ws.Range("F3").Value = 0.4336794057244024
$ws.Range("G3").Value = 13
$ws.Range("B4").Value = 0.04114622104959433
$ws.Range("C4").Value = 0.2438785208639479
$ws.Range("D4").Value = 0.09301430872978532
$ws.Range("E4").Value = 0.3049824728239072
$ws.Range("F4").Value = 0.315631483433317
$ws.Range("G4").Value = 12
$ws.Range("B5").Value = -0.004251297071472607
$ws.Range("C5").Value = 0.3196557001756482
$ws.Range("D5").Value = 0.1664324393548997
$ws.Range("E5").Value = 0.4079613209054257
$ws.Range("F5").Value = 0.4278502102499434
$ws.Range("G5").Value = 11
$ws.Range("B6").Value = -0.03152289704198964
$ws.Range("C6").Value = 0.3137577720403564
$ws.Range("D6").Value = 0.1200567451426497
$ws.Range("E6").Value = 0.3464920563918454
$ws.Range("F6").Value = 0.3637200573468097
$ws.Range("G6").Value = 10
$ws.Range("B7").Value = -0.05669611638992965
$ws.Range("C7").Value = 0.309241221828214
$ws.Range("D7").Value = 0.12362493933599
$ws.Range("E7").Value = 0.3516033835673229
$ws.Range("F7").Value = 0.3680513563860018
$ws.Range("G7").Value = 9
$ws.Range("B8").Value = -0.04410238231315148
$ws.Range("C8").Value = 0.2487243691024244
$ws.Range("D8").Value = 0.1063280554679864
$ws.Range("E8").Value = 0.3260798299005727
$ws.Range("F8").Value = 0.3453909343541705
$ws.Range("G8").Value = 8
$ws.Range("B9").Value = -0.07289500678705241
$ws.Range("C9").Value = 0.3460499651212735
$ws.Range("D9").Value = 0.2127274776992719
$ws.Range("E9").Value = 0.4612238910759848
$ws.Range("F9").Value = 0.4919174337551497
$ws.Range("G9").Value = 7
$ws.Range("B10").Value = -0.05253816173664939
$ws.Range("C10").Value = 0.2093953811075043
$ws.Range("D10").Value = 0.05818608614024737
$ws.Range("E10").Value = 0.2412179225104291
$ws.Range("F10").Value = 0.2578972532655152
$ws.Range("G10").Value = 6
$ws.Range("B11").Value = -0.009352853557239749
$ws.Range("C11").Value = 0.4545107205050175
$ws.Range("D11").Value = 0.223274494992554
$ws.Range("E11").Value = 0.4725193064759936
$ws.Range("F11").Value = 0.5281891459539978
$ws.Range("G11").Value = 5